$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The "Recorded By" column (G) stores a comma-separated list of recorder
# names/emails. Daily attendance processing moves the leading "System"
# entry so it no longer sits first in the list - it is re-inserted just
# before the final entry instead (e.g. "System, a, b" -> "a, System, b";
# "System, a" -> "a, System").
$lastRow = $ws.UsedRange.Rows.Count

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)
    $val = $cell.Value2

    if ($val -eq $null) { continue }
    $val = $val.ToString()

    if ($val.StartsWith("System, ")) {
        $parts = $val -split ", "
        $rest = $parts[1..($parts.Count - 1)]

        if ($rest.Count -eq 1) {
            $newParts = @($rest[0], "System")
        } else {
            $head = $rest[0..($rest.Count - 2)]
            $tail = $rest[$rest.Count - 1]
            $newParts = $head + @("System") + @($tail)
        }

        $cell.Value = ($newParts -join ", ")
    }
}
